$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H113").Value = 15000.223
$ws.Range("I113").Value = 18357.285
$ws.Range("J113").Value = 3250.5
$ws.Range("K113").Value = 18357.285
$ws.Range("L113").Value = 3250.5
$ws.Range("M113").Value = -15103.285
$ws.Range("N113").Value = -9758.5

$ws.Range("H134").Value = 41300
$ws.Range("J134").Value = 41300
$ws.Range("L134").Value = 41300
$ws.Range("N134").Value = -51440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3485.2952
$ws.Range("I32").Value = 2282.375
$ws.Range("K32").Value = 2282.375
$ws.Range("M32").Value = -1995.375

$ws.Range("H132").Value = 1723.3636
$ws.Range("I132").Value = 1395.8
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4187.4
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1657.4
$ws.Range("N132").Value = -20057

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 39997.5
$ws.Range("I26").Value = 39997.5
$ws.Range("K26").Value = 39997.5
$ws.Range("M26").Value = -39705.5

$ws.Range("H134").Value = 8300.532999999999
$ws.Range("I134").Value = 10355.368
$ws.Range("J134").Value = 4751.273
$ws.Range("K134").Value = 31066.104
$ws.Range("L134").Value = 14253.819
$ws.Range("M134").Value = -28531.104
$ws.Range("N134").Value = -19323.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 630.875
$ws.Range("J105").Value = 250
$ws.Range("L105").Value = 250
$ws.Range("N105").Value = -3744

$ws.Range("H107").Value = 421.51852
$ws.Range("I107").Value = 314.15
$ws.Range("J107").Value = 728.2857
$ws.Range("K107").Value = 314.15
$ws.Range("L107").Value = 728.2857
$ws.Range("M107").Value = 1605.85
$ws.Range("N107").Value = -4568.2857

$ws.Range("H134").Value = 961
$ws.Range("I134").Value = 961
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2883
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -348
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 5243.273
$ws.Range("I104").Value = 1838
$ws.Range("J104").Value = 6000
$ws.Range("K104").Value = 5514
$ws.Range("L104").Value = 18000
$ws.Range("M104").Value = -2893
$ws.Range("N104").Value = -23242

$ws.Range("H123").Value = 250002500
$ws.Range("I123").Value = 250002500
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 750007500
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -750005050
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 125004360
$ws.Range("I126").Value = 166670480
$ws.Range("K126").Value = 500011440
$ws.Range("M126").Value = -500006500

$ws.Range("H131").Value = 5564090
$ws.Range("J131").Value = 9002.165000000001
$ws.Range("L131").Value = 27006.495
$ws.Range("N131").Value = -37086.495

$ws.Range("H132").Value = 1713.3636
$ws.Range("I132").Value = 1583
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 14247
$ws.Range("L132").Value = 20700
$ws.Range("M132").Value = -11717
$ws.Range("N132").Value = -25760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2414.4707
$ws.Range("I80").Value = 2247.2
$ws.Range("J80").Value = 2484.1667
$ws.Range("K80").Value = 2247.2
$ws.Range("L80").Value = 2484.1667
$ws.Range("M80").Value = -1249.2
$ws.Range("N80").Value = -4480.1667

$ws.Range("H83").Value = 2414.4707
$ws.Range("I83").Value = 2247.2
$ws.Range("J83").Value = 2484.1667
$ws.Range("K83").Value = 11236
$ws.Range("L83").Value = 12420.8335
$ws.Range("M83").Value = -6244
$ws.Range("N83").Value = -22404.8335

$ws.Range("H102").Value = 4602.6
$ws.Range("J102").Value = 2671.3333
$ws.Range("L102").Value = 2671.3333
$ws.Range("N102").Value = -5915.3333

$ws.Range("H122").Value = 1837.2667
$ws.Range("I122").Value = 1631.7391
$ws.Range("J122").Value = 2512.5715
$ws.Range("K122").Value = 4895.2173
$ws.Range("L122").Value = 7537.7145
$ws.Range("M122").Value = -2445.2173
$ws.Range("N122").Value = -12437.7145

$ws.Range("H126").Value = 79789.62
$ws.Range("I126").Value = 3314
$ws.Range("J126").Value = 202150.6
$ws.Range("K126").Value = 9942
$ws.Range("L126").Value = 606451.8
$ws.Range("M126").Value = -7472
$ws.Range("N126").Value = -611391.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 455454.53

$ws.Range("H7").Value = 6263.364
$ws.Range("I7").Value = 3400
$ws.Range("K7").Value = 3400
$ws.Range("M7").Value = -3288

$ws.Range("H96").Value = 99995
$ws.Range("J96").Value = 99995
$ws.Range("L96").Value = 99995
$ws.Range("N96").Value = -105487

$ws.Range("H122").Value = 7993.8667
$ws.Range("I122").Value = 4780.6
$ws.Range("J122").Value = 9600.5
$ws.Range("K122").Value = 14341.8
$ws.Range("L122").Value = 28801.5
$ws.Range("M122").Value = -11891.8
$ws.Range("N122").Value = -33701.5

$ws.Range("H126").Value = 6263.364
$ws.Range("I126").Value = 3400
$ws.Range("K126").Value = 10200
$ws.Range("M126").Value = -7730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17000
$ws.Range("J31").Value = 17000
$ws.Range("L31").Value = 17000
$ws.Range("N31").Value = -17696

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H107").Value = 591.3
$ws.Range("I107").Value = 454.7143
$ws.Range("J107").Value = 910
$ws.Range("K107").Value = 1364.1429
$ws.Range("L107").Value = 2730
$ws.Range("M107").Value = 555.8571000000002
$ws.Range("N107").Value = -6570

$ws.Range("H126").Value = 5459.231
$ws.Range("I126").Value = 4685.722
$ws.Range("K126").Value = 14057.166
$ws.Range("M126").Value = -11587.166

$ws.Range("H132").Value = 3172
$ws.Range("I132").Value = 3061.0588
$ws.Range("J132").Value = 3486.3333
$ws.Range("K132").Value = 9183.1764
$ws.Range("L132").Value = 10458.9999
$ws.Range("M132").Value = -6653.1764
$ws.Range("N132").Value = -15518.9999

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
